$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right (B) 5 -> 4, Wrong (C) -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right (B) 50 -> 40, Wrong (C) -4 -> -8, Max (E) "50 / 140" -> "32 / 112"
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "32 / 112"
